$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark that sits right after the
#    title "Alteração de Característica".
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Fix the typo in the description cell: "quilometragem" (correct
#    spelling) -> "kilometragem" (the misspelling the author typed,
#    hence the later spell-check marks around it).
# ------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$descCell = $tbl.Cell(5, 2)
$cellStart = $descCell.Range.Start

# "quilometragem" starts at offset 127 inside the cell text.
$quStart = $cellStart + 127
$quEnd = $quStart + 2
$quRange = $d.Range($quStart, $quEnd)
$quRange.Text = "k"

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark at the new cursor location
#    (right after "pont", before "os."), matching where the author's
#    last edit landed.
# ------------------------------------------------------------------
$goBackPos = $cellStart + 156
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
